# Auto-generated Excel COM-interop edit script
# Applies numeric "F" column corrections across sheets 1-4,
# plus one G6 change on sheet 2 (80 -> "不可售", a type change to text).

$wb = $excel.ActiveWorkbook

# --- Worksheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 398
$ws.Range("F5").Value = 1714
$ws.Range("F6").Value = 726
$ws.Range("F7").Value = 2762
$ws.Range("F8").Value = 2128
$ws.Range("F9").Value = 880
$ws.Range("F10").Value = 2359
$ws.Range("F12").Value = 6846
$ws.Range("F13").Value = 141
$ws.Range("F15").Value = 156
$ws.Range("F16").Value = 1563
$ws.Range("F17").Value = 1358
$ws.Range("F18").Value = 1229
$ws.Range("F19").Value = 108
$ws.Range("F20").Value = 2806
$ws.Range("F21").Value = 2521
$ws.Range("F22").Value = 2521
$ws.Range("F23").Value = 817
$ws.Range("F24").Value = 1134
$ws.Range("F25").Value = 273
$ws.Range("F26").Value = 5479
$ws.Range("F27").Value = 303
$ws.Range("F30").Value = 3823
$ws.Range("F33").Value = 1729
$ws.Range("F34").Value = 1090
$ws.Range("F38").Value = 290
$ws.Range("F40").Value = 428
$ws.Range("F41").Value = 1779
$ws.Range("F42").Value = 56
$ws.Range("F43").Value = 52
$ws.Range("F45").Value = 930
$ws.Range("F46").Value = 525
$ws.Range("F50").Value = 97

# --- Worksheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 79
$ws.Range("F8").Value = 499
$ws.Range("F10").Value = 409
$ws.Range("F14").Value = 968
$ws.Range("F21").Value = 282
$ws.Range("F22").Value = 364
$ws.Range("F25").Value = 89
$ws.Range("F29").Value = 86
$ws.Range("F39").Value = 2

# --- Worksheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 3306
$ws.Range("F5").Value = 415
$ws.Range("F7").Value = 1487
$ws.Range("F9").Value = 420
$ws.Range("F10").Value = 2885
$ws.Range("F11").Value = 342
$ws.Range("F12").Value = 637
$ws.Range("F13").Value = 766
$ws.Range("F14").Value = 1271

# --- Worksheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 415
$ws.Range("F3").Value = 1487
$ws.Range("F5").Value = 398
$ws.Range("F6").Value = 420
$ws.Range("F7").Value = 2885
$ws.Range("F8").Value = 1714
$ws.Range("F9").Value = 726
$ws.Range("F10").Value = 2762
$ws.Range("F11").Value = 342
$ws.Range("F12").Value = 2128
$ws.Range("F13").Value = 880
$ws.Range("F14").Value = 2359
$ws.Range("F16").Value = 6846
$ws.Range("F17").Value = 141
$ws.Range("F18").Value = 637
$ws.Range("F19").Value = 766
$ws.Range("F20").Value = 1563
$ws.Range("F21").Value = 1358
$ws.Range("F22").Value = 1229
$ws.Range("F23").Value = 108
$ws.Range("F24").Value = 1271
$ws.Range("F25").Value = 2806
$ws.Range("F26").Value = 2521
$ws.Range("F27").Value = 364
$ws.Range("F28").Value = 817
$ws.Range("F29").Value = 1134
$ws.Range("F30").Value = 273
$ws.Range("F31").Value = 5479
$ws.Range("F32").Value = 303
$ws.Range("F33").Value = 3823
$ws.Range("F36").Value = 1729
$ws.Range("F37").Value = 1090
$ws.Range("F38").Value = 86
$ws.Range("F40").Value = 290
$ws.Range("F42").Value = 428
$ws.Range("F43").Value = 1779
$ws.Range("F44").Value = 56
$ws.Range("F46").Value = 930
$ws.Range("F47").Value = 525
$ws.Range("F51").Value = 97

# --- Special case: sheet 2, G6 becomes non-numeric text "不可售" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G6").Value = "不可售"
